$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.363.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.871.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7121'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.56%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07828'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.45%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3066'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +9.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08186'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.880.58'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.246'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.44%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7222'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.381.04'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.817'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.76%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007839'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("E20").Value = '  -0.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.131.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.47%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9992'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.767'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1471'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.965'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.931'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.363'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.10%  '

$ws.Range("E31").Value = '  -0.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.308'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.063'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05225'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.192'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7198'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.004'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.673'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01856'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.51%  '

$ws.Range("E40").Value = '  -0.93%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.172.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9139'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.991'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4295'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '71.53'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9995'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5325'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.76%  '

$ws.Range("E49").Value = '  -0.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.924'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.224'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.72%  '
